$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($row, $name, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 08:22"

# Update countries & provincias Spain data pulled from the latest refresh.
# Some countries overtook neighbours in total-case ranking, so their rows
# swapped places (country name + stats move together to keep the sheet
# sorted by total cases).

Set-RowData 21 "Israel" 9404 156 801 8532 147 6 71
Set-RowData 32 "Rumania" 4417 0 460 3754 274 6 203
Set-RowData 67 "Lituania" 912 32 8 889 11 0 15
Set-RowData 68 "Hungria" 895 78 94 743 17 11 58
Set-RowData 78 "Tunez" 623 0 25 575 67 0 23
Set-RowData 81 "Bulgaria" 581 4 42 516 21 0 23
Set-RowData 85 "Uzbekistan" 534 14 30 502 8 0 2
Set-RowData 93 "Taiwan" 379 3 67 307 0 0 5
Set-RowData 103 "Kirguistan" 270 42 33 233 5 0 4
Set-RowData 104 "Mauricio" 268 0 8 253 3 0 7
Set-RowData 105 "Estado de Palestina" 261 0 42 218 0 0 1
Set-RowData 106 "Nigeria" 254 0 44 204 2 0 6
Set-RowData 107 "Vietnam" 251 0 122 129 8 0 0
Set-RowData 108 "Montenegro" 241 0 4 235 7 0 2
Set-RowData 109 "Senegal" 237 0 105 130 1 0 2
Set-RowData 111 "Georgia" 208 12 46 159 6 0 3
Set-RowData 129 "El Salvador" 93 15 9 79 2 1 5
Set-RowData 130 "Republica de Yibuti" 90 0 9 81 0 0 0
Set-RowData 131 "Madagascar" 88 0 7 81 6 0 0
Set-RowData 132 "Guatemala" 80 3 17 60 3 0 3
Set-RowData 133 "Monaco" 79 0 4 74 4 0 1
Set-RowData 144 "Congo" 45 0 2 38 0 0 5
Set-RowData 145 "Islas Caimanes" 45 0 6 38 0 0 1
Set-RowData 161 "Birmania" 22 0 0 19 0 2 3
Set-RowData 171 "Mongolia" 16 1 4 12 0 0 0
Set-RowData 172 "Fiyi" 15 0 0 15 0 0 0
Set-RowData 173 "Dominica" 15 0 1 14 0 0 0
Set-RowData 191 "Malaui" 8 0 0 7 1 0 1
Set-RowData 192 "Islas Turcas y Caicos" 8 0 0 7 0 0 1
Set-RowData 193 "San Vicente y las Granadinas" 8 0 1 7 0 0 0
Set-RowData 207 "Burundi" 3 0 0 3 0 0 0
Set-RowData 209 "Islas Virgenes Britanicas" 3 0 0 3 0 0 0
